# Auto-generated script to apply cell value corrections described in the commit diff.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) has its own
# Leve profit calculation table; this script updates the recalculated
# price / profit columns (H-N) on the specific rows that changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 548.7273
$ws.Range("I33").Value = 583.6
$ws.Range("K33").Value = 583.6
$ws.Range("M33").Value = -354.6

$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -31240

$ws.Range("H132").Value = 6674.162
$ws.Range("I132").Value = 8348.817999999999
$ws.Range("J132").Value = 4218
$ws.Range("K132").Value = 25046.454
$ws.Range("L132").Value = 12654
$ws.Range("M132").Value = -22516.454
$ws.Range("N132").Value = -17714

$ws.Range("H138").Value = 2862.51
$ws.Range("I138").Value = 1329.8422
$ws.Range("J138").Value = 3801.8872
$ws.Range("K138").Value = 3989.5266
$ws.Range("L138").Value = 11405.6616
$ws.Range("M138").Value = 1150.4734
$ws.Range("N138").Value = -21685.6616

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23761.785
$ws.Range("I32").Value = 12233.768
$ws.Range("J32").Value = 46817.82
$ws.Range("K32").Value = 12233.768
$ws.Range("L32").Value = 46817.82
$ws.Range("M32").Value = -11946.768
$ws.Range("N32").Value = -47391.82

$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 1000
$ws.Range("M63").Value = -314

$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 5000
$ws.Range("M66").Value = -1568

$ws.Range("H74").Value = 1614.7333
$ws.Range("I74").Value = 1384.9166
$ws.Range("J74").Value = 2534
$ws.Range("K74").Value = 1384.9166
$ws.Range("L74").Value = 2534
$ws.Range("M74").Value = -510.9166
$ws.Range("N74").Value = -4282

$ws.Range("H77").Value = 1614.7333
$ws.Range("I77").Value = 1384.9166
$ws.Range("J77").Value = 2534
$ws.Range("K77").Value = 6924.583000000001
$ws.Range("L77").Value = 12670
$ws.Range("M77").Value = -2556.583000000001
$ws.Range("N77").Value = -21406

$ws.Range("H102").Value = 902.8570999999999
$ws.Range("I102").Value = 902.8570999999999
$ws.Range("K102").Value = 902.8570999999999
$ws.Range("M102").Value = 719.1429000000001

$ws.Range("H122").Value = 2273.742
$ws.Range("I122").Value = 1997.4902
$ws.Range("K122").Value = 5992.4706
$ws.Range("M122").Value = -3542.4706

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 948.875
$ws.Range("I22").Value = 747.25
$ws.Range("J22").Value = 1150.5
$ws.Range("K22").Value = 747.25
$ws.Range("L22").Value = 1150.5
$ws.Range("M22").Value = -574.25
$ws.Range("N22").Value = -1496.5

$ws.Range("H94").Value = 5902.45
$ws.Range("I94").Value = 914.1212
$ws.Range("K94").Value = 914.1212
$ws.Range("M94").Value = -463.1212

$ws.Range("H99").Value = 1665.8837
$ws.Range("I99").Value = 1313.8387
$ws.Range("J99").Value = 2575.3333
$ws.Range("K99").Value = 1313.8387
$ws.Range("L99").Value = 2575.3333
$ws.Range("M99").Value = 184.1613
$ws.Range("N99").Value = -5571.3333

$ws.Range("H105").Value = 1938.75
$ws.Range("I105").Value = 1513.8462
$ws.Range("K105").Value = 1513.8462
$ws.Range("M105").Value = 233.1538

$ws.Range("H134").Value = 1209.1794
$ws.Range("I134").Value = 1047.6666
$ws.Range("J134").Value = 1747.5555
$ws.Range("K134").Value = 3142.9998
$ws.Range("L134").Value = 5242.666499999999
$ws.Range("M134").Value = -607.9998000000001
$ws.Range("N134").Value = -10312.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2542.0894
$ws.Range("I31").Value = 1181.2222
$ws.Range("J31").Value = 3809.1035
$ws.Range("K31").Value = 1181.2222
$ws.Range("L31").Value = 3809.1035
$ws.Range("M31").Value = -886.2221999999999
$ws.Range("N31").Value = -4399.1035

$ws.Range("H34").Value = 2542.0894
$ws.Range("I34").Value = 1181.2222
$ws.Range("J34").Value = 3809.1035
$ws.Range("K34").Value = 1181.2222
$ws.Range("L34").Value = 3809.1035
$ws.Range("M34").Value = -979.2221999999999
$ws.Range("N34").Value = -4213.1035

$ws.Range("H86").Value = 230454.4
$ws.Range("I86").Value = 359435.5
$ws.Range("J86").Value = 4737.5
$ws.Range("K86").Value = 359435.5
$ws.Range("L86").Value = 4737.5
$ws.Range("M86").Value = -358312.5
$ws.Range("N86").Value = -6983.5

$ws.Range("H89").Value = 230454.4
$ws.Range("I89").Value = 359435.5
$ws.Range("J89").Value = 4737.5
$ws.Range("K89").Value = 1797177.5
$ws.Range("L89").Value = 23687.5
$ws.Range("M89").Value = -1791561.5
$ws.Range("N89").Value = -34919.5

$ws.Range("H107").Value = 1145.1111
$ws.Range("I107").Value = 934.13336
$ws.Range("J107").Value = 2200
$ws.Range("K107").Value = 934.13336
$ws.Range("L107").Value = 2200
$ws.Range("M107").Value = 985.86664
$ws.Range("N107").Value = -6040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5907.154
$ws.Range("I137").Value = 950.25
$ws.Range("J137").Value = 13838.2
$ws.Range("K137").Value = 2850.75
$ws.Range("L137").Value = 41514.60000000001
$ws.Range("M137").Value = 2249.25
$ws.Range("N137").Value = -51714.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1197.5
$ws.Range("I102").Value = 1273.5385
$ws.Range("J102").Value = 999.8
$ws.Range("K102").Value = 1273.5385
$ws.Range("L102").Value = 999.8
$ws.Range("M102").Value = 348.4614999999999
$ws.Range("N102").Value = -4243.8

$ws.Range("H132").Value = 3042.449
$ws.Range("I132").Value = 2993.2424
$ws.Range("J132").Value = 3143.9375
$ws.Range("K132").Value = 8979.727200000001
$ws.Range("L132").Value = 9431.8125
$ws.Range("M132").Value = -6449.727200000001
$ws.Range("N132").Value = -14491.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 926288.3
$ws.Range("J2").Value = 586812.3
$ws.Range("L2").Value = 586812.3
$ws.Range("N2").Value = -587036.3

$ws.Range("H132").Value = 7251840
$ws.Range("I132").Value = 10210469
$ws.Range("J132").Value = 3198.9
$ws.Range("K132").Value = 30631407
$ws.Range("L132").Value = 9596.700000000001
$ws.Range("M132").Value = -30628877
$ws.Range("N132").Value = -14656.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 986.82355
$ws.Range("I107").Value = 986.82355
$ws.Range("K107").Value = 2960.47065
$ws.Range("M107").Value = -1040.47065

$ws.Range("H122").Value = 724.8
$ws.Range("I122").Value = 716.6429000000001
$ws.Range("J122").Value = 743.8333
$ws.Range("K122").Value = 2149.9287
$ws.Range("L122").Value = 2231.4999
$ws.Range("M122").Value = 300.0712999999996
$ws.Range("N122").Value = -7131.4999
